$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.693.59"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "1.597.50"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "'211.51"
$ws.Range("E5").Value = "  +0.35%  "
$ws.Range("E6").Value = "  +0.37%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("E8").Value = "  +0.69%  "
$ws.Range("E9").Value = "  +0.94%  "
$ws.Range("D10").Value = "'19.48"
$ws.Range("E10").Value = "  +0.51%  "
$ws.Range("D11").Value = "'0.0841"
$ws.Range("E11").Value = "  +0.06%  "
$ws.Range("D12").Value = "1.820.70"
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("D13").Value = "1.631.42"
$ws.Range("E13").Value = "  +2.26%  "
$ws.Range("D14").Value = "'4.03"
$ws.Range("E14").Value = "  +0.47%  "
$ws.Range("D15").Value = "'0.523"
$ws.Range("E15").Value = "  +0.85%  "
$ws.Range("D16").Value = "'65.17"
$ws.Range("E16").Value = "  +1.34%  "
$ws.Range("D17").Value = "0.0₃0768"
$ws.Range("E17").Value = "  +5.80%  "
$ws.Range("D18").Value = "26.648.12"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("D20").Value = "'209.29"
$ws.Range("E20").Value = "  +1.14%  "
$ws.Range("D21").Value = "'7.07"
$ws.Range("E21").Value = "  +4.55%  "
$ws.Range("D22").Value = "'4.28"
$ws.Range("E22").Value = "  +1.46%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").Value = "'8.94"
$ws.Range("E24").Value = "  +1.36%  "
$ws.Range("D25").Value = "'142.95"
$ws.Range("E25").Value = "  -1.75%  "
$ws.Range("D26").Value = "'1.01"
$ws.Range("E26").Value = "  +0.45%  "
$ws.Range("E27").Value = "  -0.59%  "
$ws.Range("D28").Value = "'0.115"
$ws.Range("E28").Value = "  +0.26%  "
$ws.Range("D29").Value = "'15.37"
$ws.Range("E29").Value = "  +1.14%  "
$ws.Range("D30").Value = "'0.0519"
$ws.Range("E30").Value = "  +3.02%  "
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("E32").Value = "  +0.88%  "
$ws.Range("E33").Value = "  +1.70%  "
$ws.Range("D34").Value = "1.283.07"
$ws.Range("E34").Value = "  +0.36%  "
$ws.Range("D35").Value = "'0.616"
$ws.Range("E35").Value = "  -7.20%  "
$ws.Range("E36").Value = "  -0.41%  "
$ws.Range("D37").Value = "'1.49"
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("D38").Value = "'0.0172"
$ws.Range("E38").Value = "  +0.27%  "
$ws.Range("D39").Value = "'1.07"
$ws.Range("E39").Value = "  +17.58%  "
$ws.Range("D40").Value = "'0.825"
$ws.Range("E40").Value = "  -1.14%  "
$ws.Range("E41").Value = "  +0.81%  "
$ws.Range("E42").Value = "  -0.56%  "
$ws.Range("D43").Value = "'0.783"
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("D44").Value = "'63.28"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").Value = "1.733.30"
$ws.Range("E45").Value = "  +0.19%  "
$ws.Range("D46").Value = "'91.08"
$ws.Range("E46").Value = "  +1.34%  "
$ws.Range("D47").Value = "'1.56"
$ws.Range("E47").Value = "  -1.97%  "
$ws.Range("E48").Value = "  +0.35%  "
$ws.Range("E49").Value = "  +0.62%  "
$ws.Range("E50").Value = "  +0.22%  "
$ws.Range("D51").Value = "'7.29"
$ws.Range("E51").Value = "  -1.92%  "
